$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)
$ws.Activate()

# Row 27: clear the "Terminee" note (F27)
$ws.Range("F27").ClearContents()

# Row 28: "Faire le bouton retour" -> "Tester"
$ws.Range("F28").Value = "Tester"

# Rows 29-35: E column 0.5 -> 0.7, F column gets "Placer dans Unity"
$rows = 29,30,31,32,33,34,35
foreach ($r in $rows) {
    $ws.Range("E$r").Value = 0.7
    $ws.Range("F$r").Value = "Placer dans Unity"
}

# Update the view: scroll/selection moved from D25/F36 to A25/F29
$ws.Range("A25").Select()
$win = $excel.ActiveWindow
$win.ScrollRow = 25
$win.ScrollColumn = 1
$ws.Range("F29").Select()
